$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 10.41951533333333
$ws.Range("H2").Value = 31.258546
$ws.Range("I2").Value = 0.2236054036880828
$ws.Range("J2").Value = 0.2236054036880827
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1988276666666666
$ws.Range("N2").Value = 0.5964829999999999
$ws.Range("O2").Value = 0.08051620932651297
$ws.Range("P2").Value = 0.08051620932651297
$ws.Range("Q2").Value = 2.071687921524222
$ws.Range("R2").Value = 18.645191293718
$ws.Range("S2").Value = 0.01800385948988911
$ws.Range("T2").Value = 0.01800385948988911
$ws.Range("G3").Value = 10.41951533333333
$ws.Range("H3").Value = 31.258546
$ws.Range("I3").Value = 0.2236054036880828
$ws.Range("J3").Value = 0.2236054036880827
$ws.Range("O3").Value = 0.6640070138163814
$ws.Range("P3").Value = 0.6640070138163814
$ws.Range("Q3").Value = 17.08494875550223
$ws.Range("R3").Value = 153.76453879952
$ws.Range("S3").Value = 0.1484755563761303
$ws.Range("T3").Value = 0.1484755563761303
$ws.Range("G4").Value = 10.41951533333333
$ws.Range("H4").Value = 31.258546
$ws.Range("I4").Value = 0.2236054036880828
$ws.Range("J4").Value = 0.2236054036880827
$ws.Range("M4").Value = 0.6308773333333333
$ws.Range("N4").Value = 1.892632
$ws.Range("O4").Value = 0.2554767768571056
$ws.Range("P4").Value = 0.2554767768571056
$ws.Range("Q4").Value = 6.573436048119111
$ws.Range("R4").Value = 59.160924433072
$ws.Range("S4").Value = 0.05712598782206334
$ws.Range("T4").Value = 0.05712598782206334
$ws.Range("I5").Value = 0.3072778167205806
$ws.Range("J5").Value = 0.3072778167205806
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1988276666666666
$ws.Range("N5").Value = 0.5964829999999999
$ws.Range("O5").Value = 0.08051620932651297
$ws.Range("P5").Value = 0.08051620932651297
$ws.Range("Q5").Value = 2.846906787370666
$ws.Range("R5").Value = 25.62216108633599
$ws.Range("S5").Value = 0.02474084501246816
$ws.Range("T5").Value = 0.02474084501246816
$ws.Range("I6").Value = 0.3072778167205806
$ws.Range("J6").Value = 0.3072778167205806
$ws.Range("O6").Value = 0.6640070138163814
$ws.Range("P6").Value = 0.6640070138163814
$ws.Range("S6").Value = 0.2040346254926501
$ws.Range("T6").Value = 0.2040346254926501
$ws.Range("I7").Value = 0.3072778167205806
$ws.Range("J7").Value = 0.3072778167205806
$ws.Range("M7").Value = 0.6308773333333333
$ws.Range("N7").Value = 1.892632
$ws.Range("O7").Value = 0.2554767768571056
$ws.Range("P7").Value = 0.2554767768571056
$ws.Range("Q7").Value = 9.033194385749331
$ws.Range("R7").Value = 81.29874947174399
$ws.Range("S7").Value = 0.07850234621546238
$ws.Range("T7").Value = 0.07850234621546238
$ws.Range("G8").Value = 11.62895333333333
$ws.Range("H8").Value = 34.88686
$ws.Range("I8").Value = 0.2495602454992508
$ws.Range("J8").Value = 0.2495602454992508
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.1988276666666666
$ws.Range("N8").Value = 0.5964829999999999
$ws.Range("O8").Value = 0.08051620932651297
$ws.Range("P8").Value = 0.08051620932651297
$ws.Range("Q8").Value = 2.312157657042222
$ws.Range("R8").Value = 20.80941891338
$ws.Range("S8").Value = 0.02009364496619365
$ws.Range("T8").Value = 0.02009364496619365
$ws.Range("G9").Value = 11.62895333333333
$ws.Range("H9").Value = 34.88686
$ws.Range("I9").Value = 0.2495602454992508
$ws.Range("J9").Value = 0.2495602454992508
$ws.Range("O9").Value = 0.6640070138163814
$ws.Range("P9").Value = 0.6640070138163814
$ws.Range("Q9").Value = 19.06807230702222
$ws.Range("R9").Value = 171.6126507632
$ws.Range("S9").Value = 0.1657097533812406
$ws.Range("T9").Value = 0.1657097533812406
$ws.Range("G10").Value = 11.62895333333333
$ws.Range("H10").Value = 34.88686
$ws.Range("I10").Value = 0.2495602454992508
$ws.Range("J10").Value = 0.2495602454992508
$ws.Range("M10").Value = 0.6308773333333333
$ws.Range("N10").Value = 1.892632
$ws.Range("O10").Value = 0.2554767768571056
$ws.Range("P10").Value = 0.2554767768571056
$ws.Range("Q10").Value = 7.336443068391111
$ws.Range("R10").Value = 66.02798761551999
$ws.Range("S10").Value = 0.06375684715181661
$ws.Range("T10").Value = 0.06375684715181661
$ws.Range("G11").Value = 10.230847
$ws.Range("H11").Value = 30.692541
$ws.Range("I11").Value = 0.2195565340920857
$ws.Range("J11").Value = 0.2195565340920857
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.1988276666666666
$ws.Range("N11").Value = 0.5964829999999999
$ws.Range("O11").Value = 0.08051620932651297
$ws.Range("P11").Value = 0.08051620932651297
$ws.Range("Q11").Value = 2.034175437033666
$ws.Range("R11").Value = 18.307578933303
$ws.Range("S11").Value = 0.01767785985796206
$ws.Range("T11").Value = 0.01767785985796206
$ws.Range("G12").Value = 10.230847
$ws.Range("H12").Value = 30.692541
$ws.Range("I12").Value = 0.2195565340920857
$ws.Range("J12").Value = 0.2195565340920857
$ws.Range("O12").Value = 0.6640070138163814
$ws.Range("P12").Value = 0.6640070138163814
$ws.Range("Q12").Value = 16.77558803154667
$ws.Range("R12").Value = 150.98029228392
$ws.Range("S12").Value = 0.1457870785663604
$ws.Range("T12").Value = 0.1457870785663604
$ws.Range("G13").Value = 10.230847
$ws.Range("H13").Value = 30.692541
$ws.Range("I13").Value = 0.2195565340920857
$ws.Range("J13").Value = 0.2195565340920857
$ws.Range("M13").Value = 0.6308773333333333
$ws.Range("N13").Value = 1.892632
$ws.Range("O13").Value = 0.2554767768571056
$ws.Range("P13").Value = 0.2554767768571056
$ws.Range("Q13").Value = 6.454409473101332
$ws.Range("R13").Value = 58.089685257912
$ws.Range("S13").Value = 0.05609159566776329
$ws.Range("T13").Value = 0.05609159566776329
